$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (header row 1 untouched)
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics
$data = @(
    @("ECs",  "Col3a1", "Ddr2", "ECs",  2, 0.6666666666666666, 6.003292333333334, 18.009877, 0.002628468435338583, 0.002628468435338583, 2, 0.6666666666666666, 3.048502333333333, 9.145507, 0.02767295150267144, 0.02767295150267144, 18.30105068584878, 164.709456172639, 0.00007273747953742728, 0.00007273747953742727),
    @("ECs",  "Col3a1", "Ddr2", "FAPs", 2, 0.6666666666666666, 6.003292333333334, 18.009877, 0.002628468435338583, 0.002628468435338583, 3, 1, 83.828771, 251.486313, 0.7609603866942148, 0.7609603866942148, 503.2486182570558, 4529.237564313502, 0.002000160356968786, 0.002000160356968786),
    @("ECs",  "Col3a1", "Ddr2", "sCs",  2, 0.6666666666666666, 6.003292333333334, 18.009877, 0.002628468435338583, 0.002628468435338583, 3, 1, 23.28453333333333, 69.8536, 0.2113666618031137, 0.2113666618031138, 139.7838604452445, 1258.0547440072, 0.0005555705988323698, 0.0005555705988323698),
    @("FAPs", "Col3a1", "Ddr2", "ECs",  3, 1, 2169.811686, 6509.435058, 0.9500256210433516, 0.9500256210433515, 2, 0.6666666666666666, 3.048502333333333, 9.145507, 0.02767295150267144, 0.02767295150267144, 6614.675987664934, 59532.08388898441, 0.02629001293742798, 0.02629001293742798),
    @("FAPs", "Col3a1", "Ddr2", "FAPs", 3, 1, 2169.811686, 6509.435058, 0.9500256210433516, 0.9500256210433515, 3, 1, 83.828771, 251.486313, 0.7609603866942148, 0.7609603866942148, 181892.6469388179, 1637033.822449361, 0.7229318639585605, 0.7229318639585604),
    @("FAPs", "Col3a1", "Ddr2", "sCs",  3, 1, 2169.811686, 6509.435058, 0.9500256210433516, 0.9500256210433515, 3, 1, 23.28453333333333, 69.8536, 0.2113666618031137, 0.2113666618031138, 50523.05252972319, 454707.4727675088, 0.2008037441473632, 0.2008037441473632),
    @("sCs",  "Col3a1", "Ddr2", "ECs",  3, 1, 108.1357256666667, 324.407177, 0.04734591052130986, 0.04734591052130985, 2, 0.6666666666666666, 3.048502333333333, 9.145507, 0.02767295150267144, 0.02767295150267144, 329.6520120115266, 2966.868108103739, 0.001310201085706029, 0.001310201085706029),
    @("sCs",  "Col3a1", "Ddr2", "FAPs", 3, 1, 108.1357256666667, 324.407177, 0.04734591052130986, 0.04734591052130985, 3, 1, 83.828771, 251.486313, 0.7609603866942148, 0.7609603866942148, 9064.884983829823, 81583.9648544684, 0.03602836237868565, 0.03602836237868564),
    @("sCs",  "Col3a1", "Ddr2", "sCs",  3, 1, 108.1357256666667, 324.407177, 0.04734591052130986, 0.04734591052130985, 3, 1, 23.28453333333333, 69.8536, 0.2113666618031137, 0.2113666618031138, 2517.889908809689, 22661.0091792872, 0.01000734705691819, 0.01000734705691819)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
